{"js": "// Replace the 25 \"three-digit \u00f7 one-digit\" answer cells in the table with\n// their updated values. Each old string is unique in the document, so a\n// body-wide search + full-text replace on the single matching run is safe.\nconst replacements = [\n  [\"304\u00f74=76, 0\", \"314\u00f74=78, 2\"],\n  [\"844\u00f74=211, 0\", \"256\u00f76=42, 4\"],\n  [\"735\u00f77=105, 0\", \"452\u00f72=226, 0\"],\n  [\"569\u00f77=81, 2\", \"292\u00f75=58, 2\"],\n  [\"760\u00f74=190, 0\", \"258\u00f74=64, 2\"],\n  [\"499\u00f78=62, 3\", \"265\u00f73=88, 1\"],\n  [\"572\u00f75=114, 2\", \"155\u00f78=19, 3\"],\n  [\"749\u00f73=249, 2\", \"990\u00f75=198, 0\"],\n  [\"160\u00f73=53, 1\", \"134\u00f74=33, 2\"],\n  [\"599\u00f79=66, 5\", \"503\u00f74=125, 3\"],\n  [\"773\u00f79=85, 8\", \"458\u00f75=91, 3\"],\n  [\"628\u00f77=89, 5\", \"225\u00f74=56, 1\"],\n  [\"298\u00f79=33, 1\", \"289\u00f79=32, 1\"],\n  [\"708\u00f72=354, 0\", \"750\u00f75=150, 0\"],\n  [\"608\u00f77=86, 6\", \"714\u00f74=178, 2\"],\n  [\"519\u00f75=103, 4\", \"697\u00f73=232, 1\"],\n  [\"860\u00f75=172, 0\", \"103\u00f77=14, 5\"],\n  [\"534\u00f75=106, 4\", \"641\u00f74=160, 1\"],\n  [\"339\u00f72=169, 1\", \"591\u00f79=65, 6\"],\n  [\"480\u00f77=68, 4\", \"415\u00f74=103, 3\"],\n  [\"808\u00f79=89, 7\", \"670\u00f73=223, 1\"],\n  [\"130\u00f72=65, 0\", \"219\u00f72=109, 1\"],\n  [\"875\u00f76=145, 5\", \"417\u00f76=69, 3\"],\n  [\"666\u00f72=333, 0\", \"667\u00f73=222, 1\"],\n  [\"505\u00f77=72, 1\", \"595\u00f79=66, 1\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Pattern not found: ${oldText}`);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the 25 \"three-digit \u00f7 one-digit\" answer cells in the practice table.\n# Each \"old\" string is unique in the document body, so Find/Replace against\n# the whole document Range with MatchCase on targets exactly one run each.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @{ Old = '304\u00f74=76, 0'; New = '314\u00f74=78, 2' },\n    @{ Old = '844\u00f74=211, 0'; New = '256\u00f76=42, 4' },\n    @{ Old = '735\u00f77=105, 0'; New = '452\u00f72=226, 0' },\n    @{ Old = '569\u00f77=81, 2'; New = '292\u00f75=58, 2' },\n    @{ Old = '760\u00f74=190, 0'; New = '258\u00f74=64, 2' },\n    @{ Old = '499\u00f78=62, 3'; New = '265\u00f73=88, 1' },\n    @{ Old = '572\u00f75=114, 2'; New = '155\u00f78=19, 3' },\n    @{ Old = '749\u00f73=249, 2'; New = '990\u00f75=198, 0' },\n    @{ Old = '160\u00f73=53, 1'; New = '134\u00f74=33, 2' },\n    @{ Old = '599\u00f79=66, 5'; New = '503\u00f74=125, 3' },\n    @{ Old = '773\u00f79=85, 8'; New = '458\u00f75=91, 3' },\n    @{ Old = '628\u00f77=89, 5'; New = '225\u00f74=56, 1' },\n    @{ Old = '298\u00f79=33, 1'; New = '289\u00f79=32, 1' },\n    @{ Old = '708\u00f72=354, 0'; New = '750\u00f75=150, 0' },\n    @{ Old = '608\u00f77=86, 6'; New = '714\u00f74=178, 2' },\n    @{ Old = '519\u00f75=103, 4'; New = '697\u00f73=232, 1' },\n    @{ Old = '860\u00f75=172, 0'; New = '103\u00f77=14, 5' },\n    @{ Old = '534\u00f75=106, 4'; New = '641\u00f74=160, 1' },\n    @{ Old = '339\u00f72=169, 1'; New = '591\u00f79=65, 6' },\n    @{ Old = '480\u00f77=68, 4'; New = '415\u00f74=103, 3' },\n    @{ Old = '808\u00f79=89, 7'; New = '670\u00f73=223, 1' },\n    @{ Old = '130\u00f72=65, 0'; New = '219\u00f72=109, 1' },\n    @{ Old = '875\u00f76=145, 5'; New = '417\u00f76=69, 3' },\n    @{ Old = '666\u00f72=333, 0'; New = '667\u00f73=222, 1' },\n    @{ Old = '505\u00f77=72, 1'; New = '595\u00f79=66, 1' }\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    $found = $find.Execute($find.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n    if (-not $found) {\n        throw \"Pattern not found: $($pair.Old)\"\n    }\n}\n\n"}
